$wb = $excel.ActiveWorkbook

# Worksheet references (index is 1-based, matching sheet order in workbook.xml):
# 1 = Miguel, 2 = Profesor2, 3 = Mauricio, 4 = profesor4, 5 = Profesor5
$wsMiguel    = $wb.Worksheets.Item(1)
$wsProf2     = $wb.Worksheets.Item(2)
$wsMauricio  = $wb.Worksheets.Item(3)
$wsProf4     = $wb.Worksheets.Item(4)
$wsProf5     = $wb.Worksheets.Item(5)

# ----------------------------------------------------------------------
# Profesor2: move the "x" marks from columns D,E to columns F,G,H
# (rows 5-8 and 10-22; row 9 is the "Receso" row and stays empty)
# ----------------------------------------------------------------------
$wsProf2.Range("D5:E8").Value = ""
$wsProf2.Range("D10:E22").Value = ""
$wsProf2.Range("F5:H8").Value = "x"
$wsProf2.Range("F10:H22").Value = "x"

# ----------------------------------------------------------------------
# Mauricio: move the "x" marks from column H to columns D,E,F,G
# ----------------------------------------------------------------------
$wsMauricio.Range("H5:H8").Value = ""
$wsMauricio.Range("H10:H22").Value = ""
$wsMauricio.Range("D5:G8").Value = "x"
$wsMauricio.Range("D10:G22").Value = "x"

# ----------------------------------------------------------------------
# profesor4: move the "x" marks from column H to columns D,E,F,G
# ----------------------------------------------------------------------
$wsProf4.Range("H5:H8").Value = ""
$wsProf4.Range("H10:H22").Value = ""
$wsProf4.Range("D5:G8").Value = "x"
$wsProf4.Range("D10:G22").Value = "x"

# ----------------------------------------------------------------------
# Profesor5: move the "x" marks from column G to columns D,E,F,H
# ----------------------------------------------------------------------
$wsProf5.Range("G5:G8").Value = ""
$wsProf5.Range("G10:G22").Value = ""
$wsProf5.Range("D5:F8").Value = "x"
$wsProf5.Range("D10:F22").Value = "x"
$wsProf5.Range("H5:H8").Value = "x"
$wsProf5.Range("H10:H22").Value = "x"

# ----------------------------------------------------------------------
# Update each sheet's remembered selection (and tab-selected / active-tab
# state), finishing on Profesor5 so it ends up the active sheet/tab.
# ----------------------------------------------------------------------
$wsMiguel.Select()
$wsMiguel.Range("F9").Select()

$wsProf2.Select()
$wsProf2.Range("D5:E22").Select()

$wsMauricio.Select()
$wsMauricio.Range("H5:H22").Select()

$wsProf4.Select()
$wsProf4.Range("H5:H22").Select()

$wsProf5.Select()
$wsProf5.Range("O18").Select()
